$wb = $excel.ActiveWorkbook

# --- Sheet "Samples" ---
$ws1 = $wb.Worksheets.Item("Samples")

# New "Channels" value for the sample row, and apply a text number format
# to the whole Channels column (C) so it round-trips as text, matching the
# "abs" -> "abs600/abs700" channel-naming rework.
$ws1.Range("C2").Value = "535_485,600,700"
$ws1.Range("C1:C2").NumberFormat = "@"
$ws1.Columns.Item(3).ColumnWidth = 11.1640625

# Move the selection/active cell like the saved session did.
$ws1.Range("C3").Select()

# --- Sheet "Channel Map" ---
$ws2 = $wb.Worksheets.Item("Channel Map")

# Channel map now lists one row per discrete wavelength (600, 700) instead
# of the combined "600 700" -> "abs" mapping.
$ws2.Range("A2").Value = "600"
$ws2.Range("B2").Value = "abs600"
$ws2.Range("A3").Value = "700"
$ws2.Range("B3").Value = "abs700"
$ws2.Range("A2:A3").NumberFormat = "@"

$ws2.Range("B4").Select()

$ws1.Activate()
